$wb = $excel.ActiveWorkbook

# ----- Sheet ALC -----
$ws = $wb.Worksheets.Item(1)
if ($ws.Name -ne "ALC") { throw "Unexpected sheet name: $($ws.Name), expected ALC" }
# Row 33
$ws.Cells.Item(33,8).Value = 168.16667  # H33
$ws.Cells.Item(33,9).Value = 121.8  # I33
$ws.Cells.Item(33,11).Value = 121.8  # K33
$ws.Cells.Item(33,13).Value = 107.2  # M33
# Row 40
$ws.Cells.Item(40,8).Value = 2380.2632  # H40
$ws.Cells.Item(40,10).Value = 2740.8333  # J40
$ws.Cells.Item(40,12).Value = 2740.8333  # L40
$ws.Cells.Item(40,14).Value = -3090.8333  # N40
# Row 106
$ws.Cells.Item(106,8).Value = 1262.6666  # H106
$ws.Cells.Item(106,9).Value = 1394  # I106
$ws.Cells.Item(106,11).Value = 1394  # K106
$ws.Cells.Item(106,13).Value = -763  # M106
# Row 138
$ws.Cells.Item(138,8).Value = 3454.9453  # H138
$ws.Cells.Item(138,9).Value = 3037.7646  # I138
$ws.Cells.Item(138,10).Value = 3581.5894  # J138
$ws.Cells.Item(138,11).Value = 9113.293799999999  # K138
$ws.Cells.Item(138,12).Value = 10744.7682  # L138
$ws.Cells.Item(138,13).Value = -3973.293799999999  # M138
$ws.Cells.Item(138,14).Value = -21024.7682  # N138

# ----- Sheet ARM -----
$ws = $wb.Worksheets.Item(2)
if ($ws.Name -ne "ARM") { throw "Unexpected sheet name: $($ws.Name), expected ARM" }
# Row 4
$ws.Cells.Item(4,8).Value = 468.55554  # H4
$ws.Cells.Item(4,10).Value = 713.3333  # J4
$ws.Cells.Item(4,12).Value = 713.3333  # L4
$ws.Cells.Item(4,14).Value = -945.3333  # N4
# Row 32
$ws.Cells.Item(32,8).Value = 5089.147  # H32
$ws.Cells.Item(32,9).Value = 4657.25  # I32
$ws.Cells.Item(32,10).Value = 11999.5  # J32
$ws.Cells.Item(32,11).Value = 4657.25  # K32
$ws.Cells.Item(32,12).Value = 11999.5  # L32
$ws.Cells.Item(32,13).Value = -4370.25  # M32
$ws.Cells.Item(32,14).Value = -12573.5  # N32
# Row 45
$ws.Cells.Item(45,8).Value = 1327.1666  # H45
$ws.Cells.Item(45,9).Value = 1271.2  # I45
$ws.Cells.Item(45,11).Value = 1271.2  # K45
$ws.Cells.Item(45,13).Value = -894.2  # M45
# Row 61
$ws.Cells.Item(61,8).Value = 3865.348  # H61
$ws.Cells.Item(61,9).Value = 1660.7333  # I61
$ws.Cells.Item(61,11).Value = 1660.7333  # K61
$ws.Cells.Item(61,13).Value = -1448.7333  # M61
# Row 81
$ws.Cells.Item(81,8).Value = 0  # H81
$ws.Cells.Item(81,10).Value = 0  # J81
$ws.Cells.Item(81,12).Value = 0  # L81
$ws.Cells.Item(81,14).ClearContents()  # N81
# Row 84
$ws.Cells.Item(84,8).Value = 0  # H84
$ws.Cells.Item(84,10).Value = 0  # J84
$ws.Cells.Item(84,12).Value = 0  # L84
$ws.Cells.Item(84,14).ClearContents()  # N84
# Row 110
$ws.Cells.Item(110,8).Value = 804.8570999999999  # H110
$ws.Cells.Item(110,9).Value = 809.0769  # I110
$ws.Cells.Item(110,11).Value = 809.0769  # K110
$ws.Cells.Item(110,13).Value = 1235.9231  # M110
# Row 132
$ws.Cells.Item(132,8).Value = 1515.9231  # H132
$ws.Cells.Item(132,9).Value = 1464.4546  # I132
$ws.Cells.Item(132,10).Value = 1799  # J132
$ws.Cells.Item(132,11).Value = 4393.3638  # K132
$ws.Cells.Item(132,12).Value = 5397  # L132
$ws.Cells.Item(132,13).Value = -1863.3638  # M132
$ws.Cells.Item(132,14).Value = -10457  # N132
# Row 136
$ws.Cells.Item(136,8).Value = 3865.348  # H136
$ws.Cells.Item(136,9).Value = 1660.7333  # I136
$ws.Cells.Item(136,11).Value = 4982.199900000001  # K136
$ws.Cells.Item(136,13).Value = -2432.199900000001  # M136

# ----- Sheet BSM -----
$ws = $wb.Worksheets.Item(3)
if ($ws.Name -ne "BSM") { throw "Unexpected sheet name: $($ws.Name), expected BSM" }
# Row 20
$ws.Cells.Item(20,8).Value = 1899.25  # H20
$ws.Cells.Item(20,9).Value = 1899.4286  # I20
$ws.Cells.Item(20,11).Value = 1899.4286  # K20
$ws.Cells.Item(20,13).Value = -1652.4286  # M20

# ----- Sheet CRP -----
$ws = $wb.Worksheets.Item(4)
if ($ws.Name -ne "CRP") { throw "Unexpected sheet name: $($ws.Name), expected CRP" }
# Row 22
$ws.Cells.Item(22,8).Value = 20816.77  # H22
$ws.Cells.Item(22,9).Value = 1375  # I22
$ws.Cells.Item(22,10).Value = 43498.832  # J22
$ws.Cells.Item(22,11).Value = 1375  # K22
$ws.Cells.Item(22,12).Value = 43498.832  # L22
$ws.Cells.Item(22,13).Value = -1025  # M22
$ws.Cells.Item(22,14).Value = -44198.832  # N22
# Row 58
$ws.Cells.Item(58,8).Value = 2603.077  # H58
$ws.Cells.Item(58,9).Value = 2428.2222  # I58
$ws.Cells.Item(58,10).Value = 2996.5  # J58
$ws.Cells.Item(58,11).Value = 2428.2222  # K58
$ws.Cells.Item(58,12).Value = 2996.5  # L58
$ws.Cells.Item(58,13).Value = -2225.2222  # M58
$ws.Cells.Item(58,14).Value = -3402.5  # N58
# Row 68
$ws.Cells.Item(68,8).Value = 50000  # H68
$ws.Cells.Item(68,10).Value = 50000  # J68
$ws.Cells.Item(68,12).Value = 50000  # L68
$ws.Cells.Item(68,14).Value = -51498  # N68
# Row 71
$ws.Cells.Item(71,8).Value = 50000  # H71
$ws.Cells.Item(71,10).Value = 50000  # J71
$ws.Cells.Item(71,12).Value = 150000  # L71
$ws.Cells.Item(71,14).Value = -157488  # N71
# Row 132
$ws.Cells.Item(132,8).Value = 1050.7894  # H132
$ws.Cells.Item(132,9).Value = 981.5333000000001  # I132
$ws.Cells.Item(132,10).Value = 1310.5  # J132
$ws.Cells.Item(132,11).Value = 2944.5999  # K132
$ws.Cells.Item(132,12).Value = 3931.5  # L132
$ws.Cells.Item(132,13).Value = -414.5999000000002  # M132
$ws.Cells.Item(132,14).Value = -8991.5  # N132
# Row 136
$ws.Cells.Item(136,8).Value = 2603.077  # H136
$ws.Cells.Item(136,9).Value = 2428.2222  # I136
$ws.Cells.Item(136,10).Value = 2996.5  # J136
$ws.Cells.Item(136,11).Value = 7284.6666  # K136
$ws.Cells.Item(136,12).Value = 8989.5  # L136
$ws.Cells.Item(136,13).Value = -4734.6666  # M136
$ws.Cells.Item(136,14).Value = -14089.5  # N136
# Row 141
$ws.Cells.Item(141,8).Value = 288666.5  # H141
$ws.Cells.Item(141,10).Value = 288666.5  # J141
$ws.Cells.Item(141,12).Value = 288666.5  # L141
$ws.Cells.Item(141,14).Value = -299026.5  # N141

# ----- Sheet CUL -----
$ws = $wb.Worksheets.Item(5)
if ($ws.Name -ne "CUL") { throw "Unexpected sheet name: $($ws.Name), expected CUL" }
# Row 2
$ws.Cells.Item(2,8).Value = 122879.78  # H2
$ws.Cells.Item(2,10).Value = 101060.37  # J2
$ws.Cells.Item(2,12).Value = 606362.22  # L2
$ws.Cells.Item(2,14).Value = -606588.22  # N2
# Row 38
$ws.Cells.Item(38,8).Value = 240  # H38
$ws.Cells.Item(38,10).Value = 0  # J38
$ws.Cells.Item(38,12).Value = 0  # L38
$ws.Cells.Item(38,14).ClearContents()  # N38
# Row 64
$ws.Cells.Item(64,8).Value = 13779  # H64
$ws.Cells.Item(64,9).Value = 2497  # I64
$ws.Cells.Item(64,10).Value = 19420  # J64
$ws.Cells.Item(64,11).Value = 7491  # K64
$ws.Cells.Item(64,12).Value = 58260  # L64
$ws.Cells.Item(64,13).Value = -7221  # M64
$ws.Cells.Item(64,14).Value = -58800  # N64
# Row 67
$ws.Cells.Item(67,8).Value = 13779  # H67
$ws.Cells.Item(67,9).Value = 2497  # I67
$ws.Cells.Item(67,10).Value = 19420  # J67
$ws.Cells.Item(67,11).Value = 7491  # K67
$ws.Cells.Item(67,12).Value = 58260  # L67
$ws.Cells.Item(67,13).Value = -6555  # M67
$ws.Cells.Item(67,14).Value = -60132  # N67
# Row 68
$ws.Cells.Item(68,8).Value = 12505397  # H68
$ws.Cells.Item(68,9).Value = 500  # I68
$ws.Cells.Item(68,10).Value = 15631621  # J68
$ws.Cells.Item(68,11).Value = 1500  # K68
$ws.Cells.Item(68,12).Value = 46894863  # L68
$ws.Cells.Item(68,13).Value = -689  # M68
$ws.Cells.Item(68,14).Value = -46896485  # N68
# Row 71
$ws.Cells.Item(71,8).Value = 12505397  # H71
$ws.Cells.Item(71,9).Value = 500  # I71
$ws.Cells.Item(71,10).Value = 15631621  # J71
$ws.Cells.Item(71,11).Value = 4500  # K71
$ws.Cells.Item(71,12).Value = 140684589  # L71
$ws.Cells.Item(71,13).Value = -444  # M71
$ws.Cells.Item(71,14).Value = -140692701  # N71
# Row 113
$ws.Cells.Item(113,8).Value = 885.9  # H113
$ws.Cells.Item(113,10).Value = 751.9375  # J113
$ws.Cells.Item(113,12).Value = 2255.8125  # L113
$ws.Cells.Item(113,14).Value = -6595.8125  # N113

# ----- Sheet GSM -----
$ws = $wb.Worksheets.Item(6)
if ($ws.Name -ne "GSM") { throw "Unexpected sheet name: $($ws.Name), expected GSM" }
# Row 23
$ws.Cells.Item(23,8).Value = 750  # H23
$ws.Cells.Item(23,9).Value = 0  # I23
$ws.Cells.Item(23,10).Value = 750  # J23
$ws.Cells.Item(23,11).Value = 0  # K23
$ws.Cells.Item(23,12).Value = 750  # L23
$ws.Cells.Item(23,13).ClearContents()  # M23
$ws.Cells.Item(23,14).Value = -1196  # N23
# Row 63
$ws.Cells.Item(63,8).Value = 50057  # H63
$ws.Cells.Item(63,9).Value = 0  # I63
$ws.Cells.Item(63,10).Value = 50057  # J63
$ws.Cells.Item(63,11).Value = 0  # K63
$ws.Cells.Item(63,12).Value = 50057  # L63
$ws.Cells.Item(63,13).ClearContents()  # M63
$ws.Cells.Item(63,14).Value = -51429  # N63
# Row 66
$ws.Cells.Item(66,8).Value = 50057  # H66
$ws.Cells.Item(66,9).Value = 0  # I66
$ws.Cells.Item(66,10).Value = 50057  # J66
$ws.Cells.Item(66,11).Value = 0  # K66
$ws.Cells.Item(66,12).Value = 150171  # L66
$ws.Cells.Item(66,13).ClearContents()  # M66
$ws.Cells.Item(66,14).Value = -157035  # N66
# Row 70
$ws.Cells.Item(70,8).Value = 4212.75  # H70
$ws.Cells.Item(70,10).Value = 0  # J70
$ws.Cells.Item(70,12).Value = 0  # L70
$ws.Cells.Item(70,14).ClearContents()  # N70
# Row 73
$ws.Cells.Item(73,8).Value = 4212.75  # H73
$ws.Cells.Item(73,10).Value = 0  # J73
$ws.Cells.Item(73,12).Value = 0  # L73
$ws.Cells.Item(73,14).ClearContents()  # N73
# Row 97
$ws.Cells.Item(97,8).Value = 474.47058  # H97
$ws.Cells.Item(97,9).Value = 414.4  # I97
$ws.Cells.Item(97,11).Value = 414.4  # K97
$ws.Cells.Item(97,13).Value = 81.60000000000002  # M97

# ----- Sheet LTW -----
$ws = $wb.Worksheets.Item(7)
if ($ws.Name -ne "LTW") { throw "Unexpected sheet name: $($ws.Name), expected LTW" }
# Row 130
$ws.Cells.Item(130,8).Value = 14999.833  # H130
$ws.Cells.Item(130,10).Value = 14999.833  # J130
$ws.Cells.Item(130,12).Value = 14999.833  # L130
$ws.Cells.Item(130,14).Value = -25039.833  # N130

# ----- Sheet WVR -----
$ws = $wb.Worksheets.Item(8)
if ($ws.Name -ne "WVR") { throw "Unexpected sheet name: $($ws.Name), expected WVR" }
# Row 107
$ws.Cells.Item(107,8).Value = 1070.8  # H107
$ws.Cells.Item(107,9).Value = 588.75  # I107
$ws.Cells.Item(107,11).Value = 1766.25  # K107
$ws.Cells.Item(107,13).Value = 153.75  # M107
